# Insert a new weekly data record as row 515, pushing the existing
# rows 515-591 down by one (to 516-592). This mirrors a new Perejil
# (parsley) price observation being added to the consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 515..591 down to 516..592, opening up a blank row 515.
$ws.Rows(515).Insert()

$newRow = 515

$ws.Cells.Item($newRow, 1).Value  = 9
$ws.Cells.Item($newRow, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($newRow, 3).Value  = "Metropolitana"
$ws.Cells.Item($newRow, 4).Value  = 45127
$ws.Cells.Item($newRow, 5).Value  = 13
$ws.Cells.Item($newRow, 6).Value  = 100112044
$ws.Cells.Item($newRow, 7).Value  = "Perejil"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 70
$ws.Cells.Item($newRow, 11).Value = 12000
$ws.Cells.Item($newRow, 12).Value = 13000
$ws.Cells.Item($newRow, 13).Value = 12500
$ws.Cells.Item($newRow, 14).Value = "$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 4167
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
